# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect newly scraped counts, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 3391
    $ws.Range("F4").Value = 67
    $ws.Range("F5").Value = 1571
    $ws.Range("F6").Value = 58
}
